# Update "gh-pages" generated output at 456a3b4
# ------------------------------------------------------------
# Workbook has 4 sheets:
#   1 = 展览       (Exhibitions)
#   2 = 演出       (Performances)
#   3 = 本地生活   (Local life)
#   4 = 全部类型   (All types - aggregate of the above three, sorted by date)
# ------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------- Sheet 1: 展览 (F-column "想去人数" refresh) ----------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value  = 1525
$ws1.Range("F4").Value  = 802
$ws1.Range("F7").Value  = 1104
$ws1.Range("F8").Value  = 694
$ws1.Range("F9").Value  = 759
$ws1.Range("F10").Value = 1355
$ws1.Range("F12").Value = 1013
$ws1.Range("F13").Value = 21
$ws1.Range("F15").Value = 180
$ws1.Range("F16").Value = 43
$ws1.Range("F17").Value = 417
$ws1.Range("F20").Value = 531
$ws1.Range("F21").Value = 550
$ws1.Range("F22").Value = 740
$ws1.Range("F23").Value = 223
$ws1.Range("F24").Value = 165

# ---------- Sheet 2: 演出 (F-column refresh) ----------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value  = 982
$ws2.Range("F5").Value  = 240
$ws2.Range("F6").Value  = 14
$ws2.Range("F7").Value  = 133
$ws2.Range("F8").Value  = 60
$ws2.Range("F10").Value = 62

# ---------- Sheet 3: 本地生活 (F-column refresh) ----------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 186

# ---------- Sheet 4: 全部类型 (aggregate sheet) ----------
$ws4 = $wb.Worksheets.Item(4)

# Refresh the F-values for the rows that are not affected by the
# new-row insertion further below (rows 2-26 keep their position).
$ws4.Range("F3").Value  = 186
$ws4.Range("F4").Value  = 1525
$ws4.Range("F6").Value  = 802
$ws4.Range("F8").Value  = 982
$ws4.Range("F10").Value = 1104
$ws4.Range("F11").Value = 694
$ws4.Range("F12").Value = 759
$ws4.Range("F13").Value = 1355
$ws4.Range("F15").Value = 1013
$ws4.Range("F16").Value = 21
$ws4.Range("F18").Value = 180
$ws4.Range("F19").Value = 43
$ws4.Range("F20").Value = 417
$ws4.Range("F22").Value = 240
$ws4.Range("F25").Value = 14
$ws4.Range("F26").Value = 133

# A new event ("广州·HANAPOKO 2024 LIVE") is inserted as row 27; every
# subsequent row (old 27..38) shifts down by one (new 28..39).
$ws4.Rows.Item(27).Insert()

# Copy formatting/style from the row above into the freshly inserted
# (blank) row so the new row matches the sheet's existing look.
$ws4.Range("A26:I26").Copy($ws4.Range("A27:I27"))

$ws4.Range("A27").Value = 26
$ws4.Range("B27").Value = "2024-03-09"
$ws4.Range("C27").Value = "广州·HANAPOKO 2024 LIVE"
$ws4.Range("D27").Value = "海珠同创汇东一街11号（上冲南约11-2） 声音共和Livehouse"
$ws4.Range("E27").Value = "2024.03.09 14:00-03.09 15:30"
$ws4.Range("F27").Value = 133
$ws4.Range("G27").Value = 380
$ws4.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=81279"
$ws4.Range("I27").Value = "//i2.hdslb.com/bfs/openplatform/202401/tMZ1Jp2G1705992352054.jpeg"

# Rows 28-39 already hold the shifted-down content (old rows 27-38);
# only their F/G values need refreshing to the latest figures.
$ws4.Range("F28").Value = 531
$ws4.Range("F29").Value = 550
$ws4.Range("F30").Value = 740
$ws4.Range("F31").Value = 223
$ws4.Range("F32").Value = 60
$ws4.Range("F33").Value = 165
$ws4.Range("F34").Value = 581
$ws4.Range("G34").Value = "已售罄"
$ws4.Range("F35").Value = 62
$ws4.Range("F36").Value = 62
$ws4.Range("F37").Value = 12
$ws4.Range("F38").Value = 361
$ws4.Range("F39").Value = 6
